# 16Feb2022 Selenium Testng concepts
#
# 1) The "datetimeFigureOut" Date placeholder (master + every slide layout)
#    gets re-cached from 04-09-2021 to 11-02-2022.
# 2) Slide 6 ("@DataProvider and @Factory Annotations"), the body placeholder
#    paragraph that starts "In simple words, we can say that @DataProvider..."
#    has three adjacent runs merged into one run (same visible text).

$p = $ppt.ActivePresentation

$oldDate = "04-09-2021"
$newDate = "11-02-2022"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
            $shpTr = $shp.TextFrame.TextRange
            if ($shpTr.Text -eq $oldDate) {
                $shpTr.Text = $newDate
            }
        }
    }
}

# Slide master's own Date Placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every custom (slide) layout has its own cached Date Placeholder too.
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# Slide 6: merge the "... is meant to / send parameters / to the methods,
# Where the ..." runs into a single run (formatting of the first run wins).
$slide6 = $p.Slides.Item(6)
$body = $slide6.Shapes.Item(2)
$tr = $body.TextFrame.TextRange
$fullText = $tr.Text

$oldSpan = " is meant to send parameters to the methods, Where the "
$needle = " is meant to "
$startIdx0 = $fullText.IndexOf($needle)
if ($startIdx0 -ge 0) {
    $startPos = $startIdx0 + 1
    $span = $tr.Characters($startPos, $oldSpan.Length)
    if ($span.Text -eq $oldSpan) {
        $span.Text = $oldSpan
    }
}
